$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '48.276.84'
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = '  +1.15%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.508.96'
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = '  +0.04%  '
$ws.Range("E4").Value = '  -0.12%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '319.84'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  -0.82%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '107.42'
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = '  -0.48%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.526'
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = '  +0.48%  '
$ws.Range("E8").Value = '  -0.06%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.539'
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = '  -3.62%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '39.06'
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = '  -3.05%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '19.90'
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = '  +1.91%  '
$ws.Range("E12").Value = '  -0.76%  '
$ws.Range("E13").Value = '  -0.31%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '7.10'
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = '  -0.87%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '2.903.21'
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = '  +0.13%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '2.512.82'
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = '  -0.03%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.834'
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = '  -1.88%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '48.160.26'
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = '  +1.10%  '
$ws.Range("B19").Value = 'InternetComputer(DFINITY)'
$ws.Range("C19").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '13.00'
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = '  -2.43%  '
$ws.Range("B20").Value = 'ImmutableX'
$ws.Range("C20").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '2.96'
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = '  +6.49%  '
$ws.Range("E21").Value = '  +0.52%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.0₃0939'
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = '  -0.16%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '71.26'
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = '  +0.44%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '272.87'
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = '  +10.33%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.53'
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = '  -1.70%  '
$ws.Range("E26").Value = '  -0.05%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '25.96'
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = '  +0.73%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.29'
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = '  +9.92%  '
$ws.Range("E29").Value = '  +2.16%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '9.76'
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = '  -4.57%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '34.71'
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = '  -0.36%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '49.57'
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = '  -0.57%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '19.25'
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = '  -3.81%  '
$ws.Range("E34").Value = '  -0.11%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '5.31'
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = '  -1.33%  '
$ws.Range("E36").Value = '  -0.21%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.95'
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = '  -0.81%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '4.63'
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = '  -1.74%  '
$ws.Range("E39").Value = '  -3.15%  '
$ws.Range("E40").Value = '  -1.05%  '
$ws.Range("E41").Value = '  +1.33%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '120.19'
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = '  +1.20%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '21.96'
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = '  -1.27%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.0304'
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = '  +2.31%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '2.004.97'
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = '  -0.02%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '3.20'
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = '  +3.40%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.89'
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = '  +4.85%  '
$ws.Range("E48").Value = '  -0.91%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '8.98'
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = '  -1.16%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '5.24'
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = '  +1.54%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '79.42'
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = '  +2.63%  '
